$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 16 data rows (old rows 2-17), shifting all data below
# up by 16 rows. This corrects the evaluation series so it starts at the
# correct revision date (old A18/B18 becomes the new A2/B2) and fixes the
# off-by-16 typo noted in the commit message.
$ws.Range("A2:B17").EntireRow.Delete() | Out-Null
